$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 59 down to new row 60 (same formatting/styles as row 59)
$ws.Range("A59:K59").Copy()
$ws.Range("A60:K60").PasteSpecial()

# Update the "G11M2-awake" row (59) and the new row (60) to point at the
# new quick-test raw/save data locations.
$ws.Range("B60").Value2 = "G26M4"
$ws.Range("D59").Value2 = "D:\ProcessedData\TestData"
$ws.Range("D60").Value2 = "D:\ProcessedData\TestData"
$ws.Range("C59").Value2 = "C:\Users\Nischal\Documents\QuickTestData"
$ws.Range("C60").Value2 = "C:\Users\Nischal\Documents\QuickTestData"

# Fill in the rest of row 60's values
$ws.Range("A60").Value2 = 190905
$ws.Range("H60").Value2 = 100
$ws.Range("I60").Value2 = 9
$ws.Range("J60").Value2 = 9

[void]$ws.Range("D63").Select()
